$wb = $excel.ActiveWorkbook

# --- DatosCuenta sheet ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokePreCinco"
$wsCuenta.Range("B2").Value = "SmokePreNameCinco"
$wsCuenta.Range("C2").Value = 27100112
$wsCuenta.Range("D2").Value = 114

# --- DatosHogar sheet ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 634

# --- DatosMotor sheet ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMA015"
$wsMotor.Range("B2").Value = "ABC12SSMA015"
$wsMotor.Range("C2").Value = "ZAZ123SSMA015"

# --- DatosAP sheet ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200114
$wsAP.Range("E10").Select()

# Restore original active sheet/tab selection (DatosAP selection change
# should not shift which sheet tab is marked as selected)
$wsCuenta.Activate()
